# Daily attendance processing - 2025-11-21 15:21:18
#
# For every row in the "Recorded By" column (G), when the cell contains
# multiple comma-separated recorder names/emails, reverse the order of
# those names. Cells that only contain a single value are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count()

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($val -ne $null -and $val -ne "") {
        if ($val.Contains(",")) {
            $parts = $val.Split(",")

            $trimmed = @()
            foreach ($p in $parts) {
                $trimmed += $p.Trim()
            }

            # Build the reversed list manually (array index walk),
            # since [array]::Reverse() does not mutate in this runtime.
            $count = $trimmed.Count
            $reversed = @()
            for ($i = $count - 1; $i -ge 0; $i--) {
                $reversed += $trimmed[$i]
            }

            $newVal = [string]::Join(", ", $reversed)
            $cell.Value = $newVal
        }
    }
}
